# Append-style update of the "ランサーズ" sheet:
#   - refresh the "取得日時" timestamp on every existing data row
#   - insert two newly-scraped listings at their sorted (by score) positions
#   - rebuild the hyperlinks on column F so they keep pointing at the right row
#
# Matches commit: "Append: 2025-09-29 01:17 JST"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-29 01:17:37"

# --- Drop the existing hyperlinks up front -------------------------------
# Hyperlinks.Delete() (called off any single-cell Range on the sheet) clears
# every hyperlink on the worksheet; we rebuild the full, correctly-ordered
# set once the new rows are in place.
$ws.Range("A1").Hyperlinks.Delete()

# --- Make room for the two new listings -----------------------------------
# "【急募】リスト抽出ツール開発..." (score 128) slots in right after the
# existing score-158 row, i.e. becomes the new row 4.
$ws.Rows.Item(4).Insert()
# "【急募】新しい口コミサイトの構築..." (score 38) slots in right after the
# existing score-68 row. After the first insert that landing spot is row 9.
$ws.Rows.Item(9).Insert()

# --- Refresh the timestamp column for every data row -----------------------
$ws.Range("A2:A11").Value = $newTimestamp

# --- Fill in the newly inserted row 4 --------------------------------------
$ws.Range("B4").Value = "【急募】リスト抽出ツール開発のフリーランスを探しています!"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5402362"
$ws.Range("G4").Value = 128
$ws.Range("H4").Value = "◆ツール,開発"

# --- Fill in the newly inserted row 9 --------------------------------------
$ws.Range("B9").Value = "【急募】新しい口コミサイトの構築をお手伝いください!"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5402277"
$ws.Range("G9").Value = 38
$ws.Range("H9").Value = "◇サイト"

# --- Rebuild hyperlinks on column F for every data row (2..11) ------------
$urls = @{
    2  = "https://www.lancers.jp/work/detail/5217096"
    3  = "https://www.lancers.jp/work/detail/5394578"
    4  = "https://www.lancers.jp/work/detail/5402362"
    5  = "https://www.lancers.jp/work/detail/5402230"
    6  = "https://www.lancers.jp/work/detail/5402140"
    7  = "https://www.lancers.jp/work/detail/5402038"
    8  = "https://www.lancers.jp/work/detail/5402182"
    9  = "https://www.lancers.jp/work/detail/5402277"
    10 = "https://www.lancers.jp/work/detail/5399347"
    11 = "https://www.lancers.jp/work/detail/5402130"
}

foreach ($r in 2..11) {
    $ws.Hyperlinks.Add($ws.Range("F$r"), $urls[$r])
    $ws.Range("F$r").Style = "Hyperlink"
}
